$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Flat array of values for range E2:T17 (16 rows x 16 cols), row-major order
$vals = @(
    3,1,3.545238666666667,10.635716,0.0286068503401377,0.0286068503401377,1,0.3333333333333333,0.323322,0.969966,0.1201177690747821,0.1201177690747821,1.146253656184,10.316282905656,0.003436191043113514,0.003436191043113513
    3,1,3.545238666666667,10.635716,0.0286068503401377,0.0286068503401377,3,1,1.141707,3.425121,0.4241570254330428,0.4241570254330428,4.047623802404,36.428614221636,0.01213379654728104,0.01213379654728104
    3,1,3.545238666666667,10.635716,0.0286068503401377,0.0286068503401377,2,0.6666666666666666,0.6218673333333333,1.865602,0.2310307270785282,0.2310307270785282,2.204668115670222,19.842013041032,0.006609061433508655,0.006609061433508655
    3,1,3.545238666666667,10.635716,0.0286068503401377,0.0286068503401377,2,0.6666666666666666,0.604812,1.814436,0.2246944784136469,0.2246944784136469,2.144202888464,19.297825996176,0.006427801316234497,0.006427801316234497
    3,1,6.418049333333333,19.254148,0.05178781854111765,0.05178781854111765,1,0.3333333333333333,0.323322,0.969966,0.1201177690747821,0.1201177690747821,2.075096546552,18.675868918968,0.006220637228408691,0.00622063722840869
    3,1,6.418049333333333,19.254148,0.05178781854111765,0.05178781854111765,3,1,1.141707,3.425121,0.4241570254330428,0.4241570254330428,7.327531850212,65.947786651908,0.02196616706606665,0.02196616706606665
    3,1,6.418049333333333,19.254148,0.05178781854111765,0.05178781854111765,2,0.6666666666666666,0.6218673333333333,1.865602,0.2310307270785282,0.2310307270785282,3.991175224121778,35.920577017096,0.01196457737136529,0.01196457737136529
    3,1,6.418049333333333,19.254148,0.05178781854111765,0.05178781854111765,2,0.6666666666666666,0.604812,1.814436,0.2246944784136469,0.2246944784136469,3.881713253392,34.935419280528,0.01163643687527702,0.01163643687527702
    2,0.6666666666666666,0.314582,0.943746,0.002538390511847401,0.002538390511847401,1,0.3333333333333333,0.323322,0.969966,0.1201177690747821,0.1201177690747821,0.101711281404,0.915401532636,0.0003049058053237042,0.0003049058053237041
    2,0.6666666666666666,0.314582,0.943746,0.002538390511847401,0.002538390511847401,3,1,1.141707,3.425121,0.4241570254330428,0.4241570254330428,0.359160471474,3.232444243266,0.001076676168892653,0.001076676168892653
    2,0.6666666666666666,0.314582,0.943746,0.002538390511847401,0.002538390511847401,2,0.6666666666666666,0.6218673333333333,1.865602,0.2310307270785282,0.2310307270785282,0.1956282694546667,1.760654425092,0.0005864462055613424,0.0005864462055613424
    2,0.6666666666666666,0.314582,0.943746,0.002538390511847401,0.002538390511847401,2,0.6666666666666666,0.604812,1.814436,0.2246944784136469,0.2246944784136469,0.190262968584,1.712366717256,0.0005703623320697019,0.0005703623320697019
    3,1,113.65184,340.95552,0.9170669406068972,0.9170669406068972,1,0.3333333333333333,0.323322,0.969966,0.1201177690747821,0.1201177690747821,36.74614021248,330.71526191232,0.1101560349979362,0.1101560349979362
    3,1,113.65184,340.95552,0.9170669406068972,0.9170669406068972,3,1,1.141707,3.425121,0.4241570254330428,0.4241570254330428,129.75710129088,1167.81391161792,0.3889803856508025,0.3889803856508025
    3,1,113.65184,340.95552,0.9170669406068972,0.9170669406068972,2,0.6666666666666666,0.6218673333333333,1.865602,0.2310307270785282,0.2310307270785282,70.67636666922667,636.0873000230399,0.2118706420680929,0.2118706420680929
    3,1,113.65184,340.95552,0.9170669406068972,0.9170669406068972,2,0.6666666666666666,0.604812,1.814436,0.2246944784136469,0.2246944784136469,68.73799665407999,618.6419698867201,0.2060598778900656,0.2060598778900656
)

$numCols = 16
$startRow = 2
$startCol = 5  # column E

for ($k = 0; $k -lt $vals.Count; $k++) {
    $r = $startRow + [Math]::Floor($k / $numCols)
    $c = $startCol + ($k % $numCols)
    $ws.Cells.Item($r, $c).Value = $vals[$k]
}
